# Apply cell updates per the symbol-list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '242.70'

# Row 4
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '5.447'

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '0.05696'

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '3.416'

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '6.272'

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '1.121'
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '7FTXTokenFTT'

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.1421'

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07307'

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '0.03090'

# Row 14
$ws.Cells.Item(14,2).NumberFormat = "@"
$ws.Cells.Item(14,2).Value = 'BitMartToken'
$ws.Cells.Item(14,3).NumberFormat = "@"
$ws.Cells.Item(14,3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '0.09368'
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = '13BitMartTokenBMX'

# Row 15
$ws.Cells.Item(15,2).NumberFormat = "@"
$ws.Cells.Item(15,2).Value = 'MCDex'
$ws.Cells.Item(15,3).NumberFormat = "@"
$ws.Cells.Item(15,3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '3.928'
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '14MCDexMCB'

# Row 16
$ws.Cells.Item(16,2).NumberFormat = "@"
$ws.Cells.Item(16,2).Value = 'BitForexToken'
$ws.Cells.Item(16,3).NumberFormat = "@"
$ws.Cells.Item(16,3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.001582'
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = '15BitForexTokenBF'

# Row 17
$ws.Cells.Item(17,2).NumberFormat = "@"
$ws.Cells.Item(17,2).Value = 'CoinExToken'
$ws.Cells.Item(17,3).NumberFormat = "@"
$ws.Cells.Item(17,3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.04800'
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = '16CoinExTokenCET'

# Row 18
$ws.Cells.Item(18,2).NumberFormat = "@"
$ws.Cells.Item(18,2).Value = 'One'
$ws.Cells.Item(18,3).NumberFormat = "@"
$ws.Cells.Item(18,3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '0.01077'
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = '17OneONEBestin24h'

# Row 19
$ws.Cells.Item(19,2).NumberFormat = "@"
$ws.Cells.Item(19,2).Value = 'TigerCash'
$ws.Cells.Item(19,3).NumberFormat = "@"
$ws.Cells.Item(19,3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '0.006274'
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = '18TigerCashTCH'

# Row 20
$ws.Cells.Item(20,2).NumberFormat = "@"
$ws.Cells.Item(20,2).Value = 'HotbitToken'
$ws.Cells.Item(20,3).NumberFormat = "@"
$ws.Cells.Item(20,3).Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '0.004073'
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = '19HotbitTokenHTB'

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '0.0009947'

# Row 22
$ws.Cells.Item(22,2).NumberFormat = "@"
$ws.Cells.Item(22,2).Value = 'NitroEx'
$ws.Cells.Item(22,3).NumberFormat = "@"
$ws.Cells.Item(22,3).Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '0.0001499'
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = '21NitroExNTX'

# Row 23
$ws.Cells.Item(23,2).NumberFormat = "@"
$ws.Cells.Item(23,2).Value = 'LEO'
$ws.Cells.Item(23,3).NumberFormat = "@"
$ws.Cells.Item(23,3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '3.734'
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '22LEOLEO'

# Row 24
$ws.Cells.Item(24,2).NumberFormat = "@"
$ws.Cells.Item(24,2).Value = 'BTSEToken'
$ws.Cells.Item(24,3).NumberFormat = "@"
$ws.Cells.Item(24,3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '2.151'
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = '23BTSETokenBTSE'

# Row 25
$ws.Cells.Item(25,2).NumberFormat = "@"
$ws.Cells.Item(25,2).Value = 'BitpandaEcosystemToken'
$ws.Cells.Item(25,3).NumberFormat = "@"
$ws.Cells.Item(25,3).Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '0.3259'
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = '24BitpandaEcosystemTokenBEST'

# Row 26
$ws.Cells.Item(26,2).NumberFormat = "@"
$ws.Cells.Item(26,2).Value = 'ProBitToken'
$ws.Cells.Item(26,3).NumberFormat = "@"
$ws.Cells.Item(26,3).Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '0.1299'
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = '25ProBitTokenPROB'

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '0.0003996'

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.006669'

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.002998'

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.006486'
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = '43LocalTradersLCTWorstin24h'

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.00005605'

# Row 47
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.3898'

# Row 48
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = '47BOLOBOLO'

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.01009'

